# Before this edit, the deck's (only) theme part - ppt/theme/theme1.xml,
# linked from the slide master / every slide - carries the "Integral"
# design's "Red Violet" colour scheme. The commit swaps those 12 theme
# colours for the stock "Office Theme" palette (the palette that, before
# the edit, only the notes master's theme part referenced). The font
# scheme and format scheme are already identical between the old and new
# palette, so only the colour scheme slots need to change.

$p      = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# VBA-style RGB(r,g,b) -> the little-endian 0x00BBGGRR integer that
# ColorFormat.RGB reads/writes.
function RgbValue([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office" colour scheme, in clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    (RgbValue 0x00 0x00 0x00),  # 1  dk1
    (RgbValue 0xFF 0xFF 0xFF),  # 2  lt1
    (RgbValue 0x44 0x54 0x6A),  # 3  dk2
    (RgbValue 0xE7 0xE6 0xE6),  # 4  lt2
    (RgbValue 0x5B 0x9B 0xD5),  # 5  accent1
    (RgbValue 0xED 0x7D 0x31),  # 6  accent2
    (RgbValue 0xA5 0xA5 0xA5),  # 7  accent3
    (RgbValue 0xFF 0xC0 0x00),  # 8  accent4
    (RgbValue 0x44 0x72 0xC4),  # 9  accent5
    (RgbValue 0x70 0xAD 0x47),  # 10 accent6
    (RgbValue 0x05 0x63 0xC1),  # 11 hlink
    (RgbValue 0x95 0x4F 0x72)   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}

# Best-effort: also rename the colour scheme / theme to match the stock
# "Office" naming (no-op on hosts where these are read-only).
try { $themeColors.Name = "Office" } catch {}
try { $master.Theme.Name = "Office Theme" } catch {}
